$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stats columns (runs/balls/fours/sixes) are stored as text in the source
# sheet (numbers-as-text, flagged ignorable by Excel). Force the range to
# Text format first so re-assigning numeric-looking strings doesn't get
# silently reinterpreted as real numbers.
$ws.Range("C2:F12").NumberFormat = "@"

# New per-match row values (existing rows 2-11 reshuffled + one new match
# appended as row 12), per the "updated activity till excel form" edit.
$stats = @(
  @(2,  "7",  "6",  "1", "0"),
  @(3,  "31", "11", "2", "3"),
  @(4,  "21", "8",  "3", "1"),
  @(5,  "1",  "2",  "0", "0"),
  @(6,  "35", "30", "4", "0"),
  @(7,  "25", "10", "3", "1"),
  @(8,  "33", "13", "0", "4"),
  @(9,  "50", "35", "5", "2"),
  @(10, "10", "5",  "2", "0"),
  @(11, "7",  "6",  "1", "0"),
  @(12, "12", "9",  "1", "0")
)

foreach ($row in $stats) {
  $r = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
}

# Row 12 is a brand new row, so it also needs the player/team labels that
# every other row already carries.
$ws.Cells.Item(12, 1).Value = "Ravindra Jadeja "
$ws.Cells.Item(12, 2).Value = "Chennai Super Kings"
